$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Find the last used row in column A (currently row 282 -> "08-10-2021")
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

$dates = @("09-10-2021", "10-10-2021", "11-10-2021", "12-10-2021", "13-10-2021", "14-10-2021")
$firstNewRow = $lastRow + 1
$lastNewRow = $lastRow + $dates.Count

# Pre-format column A for the new rows as Text so Excel doesn't
# auto-convert the dd-mm-yyyy-looking strings into date serial numbers;
# the source data keeps these values as plain text (shared strings),
# matching the rest of column A.
for ($i = $firstNewRow; $i -le $lastNewRow; $i++) {
    $ws.Cells.Item($i, 1).NumberFormat = "@"
}

$r = $firstNewRow
foreach ($d in $dates) {
    $ws.Cells.Item($r, 1).Value = $d
    $ws.Cells.Item($r, 2).Value = 3068
    $ws.Cells.Item($r, 3).Value = 204
    $r = $r + 1
}

# Restore the default (Normal) styling so the new cells keep the same
# plain formatting as the rest of the data rows.
for ($i = $firstNewRow; $i -le $lastNewRow; $i++) {
    $ws.Cells.Item($i, 1).Style = "Normal"
}
